# "invalid user name pwd test case"
# Adds a second worksheet ("InvalidLogin") after the existing "ValidLogin"
# sheet, populates it with a UserName/Password header row plus an invalid
# admin1/manager123 login row, and updates the selection/active-sheet state
# to match (InvalidLogin becomes the active/selected tab).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new sheet right after ValidLogin.
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "InvalidLogin"

# Header row (reuses the same shared strings as ValidLogin).
$ws2.Range("A1").Value = "UserName"
$ws2.Range("B1").Value = "Password"

# Invalid credentials test data.
$ws2.Range("A2").Value = "admin1"
$ws2.Range("B2").Value = "manager123"

# Size column B to fit the longer "manager123" value.
$ws2.Columns.Item(2).AutoFit()

# Clear the old selection/active marker on ValidLogin and select A1:B2 there.
$ws1.Range("A1:B2").Select()

# Make InvalidLogin the active sheet with its own selection.
$ws2.Activate()
$ws2.Range("H12").Select()
